# Apply the "Updated symbol list" price refresh to the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Column D (Price) updates -------------------------------------------
Set-TextValue "D2"  "281.55"
Set-TextValue "D3"  "20.71"
Set-TextValue "D4"  "6.223"
Set-TextValue "D5"  "0.06156"
Set-TextValue "D6"  "3.580"
Set-TextValue "D7"  "6.573"
Set-TextValue "D8"  "1.498"
Set-TextValue "D9"  "0.8191"
Set-TextValue "D11" "0.1644"
Set-TextValue "D12" "0.08439"
Set-TextValue "D13" "0.03538"
Set-TextValue "D14" "0.03189"
Set-TextValue "D15" "0.09143"
Set-TextValue "D16" "3.707"
Set-TextValue "D17" "0.001641"
Set-TextValue "D18" "0.04716"
Set-TextValue "D19" "0.006527"
Set-TextValue "D20" "0.006165"
Set-TextValue "D21" "0.001069"
Set-TextValue "D22" "0.0001612"
Set-TextValue "D23" "3.783"
Set-TextValue "D25" "0.3356"
Set-TextValue "D40" "0.04696"
Set-TextValue "D41" "0.007193"
Set-TextValue "D44" "0.01111"
Set-TextValue "D45" "0.00006624"
Set-TextValue "D48" "0.002971"

# --- Rows 42/43 swap places (BKEXToken <-> CEJI) with refreshed prices --
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.004505"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1101"
$ws.Range("E43").Value = "42BKEXTokenBKK"
